# Scope-of-Work doc: the fiber-footage figure was corrected from
# "21348" feet to "21378" feet ("348" -> "378" in the run of text
# immediately following "Place approximately  21").
$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("21348", $false, $false, $false, $false, $false, `
                                 $true, 1, $false, "21378", 2)
